# Applies the OOXML diff to "Översikt OKÄNT" worksheet:
#  - insert two new rows at the top of the data block (rows 2-3)
#  - fill them with the two new logging cases
#  - bump the "Förändrad" (C column) date for every data row to 45247

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new rows, pushing all existing data rows down by 2 ---
$ws.Rows("2:3").Insert()
$ws.Rows("2:3").RowHeight = 15

# --- 2. Populate new row 2: A 57619-2023 ---
$ws.Range("A2").Value2 = "A 57619-2023"
$ws.Range("B2").Value2 = 45243
$ws.Range("B2").NumberFormat = "YYYY-MM-DD"
$ws.Range("C2").Value2 = 45247
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"
$ws.Range("D2").Value2 = "OKÄNT"
$ws.Range("E2").Value2 = "OKÄNT"
$ws.Range("F2").Value2 = "SCA"
$ws.Range("G2").Value2 = 17.1
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = 0
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 0
$ws.Range("Q2").Value2 = 1
$ws.Range("R2").Value2 = "Harticka"
$ws.Range("R2").WrapText = $true
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/artfynd/A 57619-2023 artfynd.xlsx", "A 57619-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/kartor/A 57619-2023 karta.png", "A 57619-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomål/A 57619-2023 FSC-klagomål.docx", "A 57619-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomålsmail/A 57619-2023 FSC-klagomål mail.docx", "A 57619-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsyn/A 57619-2023 tillsynsbegäran.docx", "A 57619-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsynsmail/A 57619-2023 tillsynsbegäran mail.docx", "A 57619-2023")'

# --- 3. Populate new row 3: A 57664-2023 ---
$ws.Range("A3").Value2 = "A 57664-2023"
$ws.Range("B3").Value2 = 45243
$ws.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws.Range("C3").Value2 = 45247
$ws.Range("C3").NumberFormat = "YYYY-MM-DD"
$ws.Range("D3").Value2 = "OKÄNT"
$ws.Range("E3").Value2 = "OKÄNT"
$ws.Range("F3").Value2 = "SCA"
$ws.Range("G3").Value2 = 6.1
$ws.Range("H3").Value2 = 0
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 1
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 1
$ws.Range("R3").Value2 = "Norsk näverlav"
$ws.Range("R3").WrapText = $true
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/artfynd/A 57664-2023 artfynd.xlsx", "A 57664-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/kartor/A 57664-2023 karta.png", "A 57664-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomål/A 57664-2023 FSC-klagomål.docx", "A 57664-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomålsmail/A 57664-2023 FSC-klagomål mail.docx", "A 57664-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsyn/A 57664-2023 tillsynsbegäran.docx", "A 57664-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsynsmail/A 57664-2023 tillsynsbegäran mail.docx", "A 57664-2023")'

# --- 4. Bump "Förändrad" date (column C) to 45247 for every remaining data row (now rows 4-28) ---
for ($r = 4; $r -le 28; $r++) {
    $ws.Range("C$r").Value2 = 45247
}
